$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-11-02 Saturday" "2024-11-03 Sunday"

Replace-Text "720÷2=360, 0" "472÷3=157, 1"
Replace-Text "679÷7=97, 0" "780÷3=260, 0"
Replace-Text "567÷9=63, 0" "767÷7=109, 4"
Replace-Text "553÷8=69, 1" "998÷5=199, 3"
Replace-Text "184÷3=61, 1" "169÷6=28, 1"

Replace-Text "994÷7=142, 0" "440÷8=55, 0"
Replace-Text "849÷8=106, 1" "821÷9=91, 2"
Replace-Text "940÷7=134, 2" "362÷7=51, 5"
Replace-Text "566÷2=283, 0" "516÷9=57, 3"
Replace-Text "708÷9=78, 6" "560÷3=186, 2"

Replace-Text "143÷9=15, 8" "861÷7=123, 0"
Replace-Text "459÷8=57, 3" "820÷7=117, 1"
Replace-Text "403÷9=44, 7" "551÷9=61, 2"
Replace-Text "797÷8=99, 5" "370÷5=74, 0"
Replace-Text "676÷3=225, 1" "672÷3=224, 0"

Replace-Text "422÷3=140, 2" "897÷7=128, 1"
Replace-Text "183÷6=30, 3" "236÷3=78, 2"
Replace-Text "514÷4=128, 2" "885÷4=221, 1"
Replace-Text "782÷6=130, 2" "595÷2=297, 1"
Replace-Text "988÷2=494, 0" "664÷3=221, 1"

Replace-Text "366÷4=91, 2" "321÷6=53, 3"
Replace-Text "707÷6=117, 5" "810÷8=101, 2"
Replace-Text "445÷5=89, 0" "418÷6=69, 4"
Replace-Text "642÷9=71, 3" "207÷5=41, 2"
Replace-Text "904÷4=226, 0" "838÷2=419, 0"
